$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column -> new value for both row 2 and row 3 (values are identical across the two rows)
$values = [ordered]@{
    "G"  = -0.5910326086956522
    "H"  = -0.5910326086956522
    "I"  = -1.546195652173913
    "J"  = -1.546195652173913
    "K"  = -326.3
    "L"  = -4.433423913043478
    "U"  = 143.8
    "V"  = 0.7420020639834881
    "W"  = -0.3937017374517375
    "X"  = 0.2372606039529525
    "Y"  = -0.6309623414046899
    "Z"  = 0.03577852316367702
    "AA" = -0.05532059695688105
    "AB" = 0.0378954162365609
    "AC" = -0.09321601319344194
    "AD" = 1589.2
    "AE" = 0
    "AF" = 1589.2
    "AG" = 1445.4
    "AH" = 0.8913067863151991
    "AI" = 0.7631578947368421
    "AJ" = 0.8817715959004393
    "AK" = 0.7455896007428041
    "AL" = 113.6
    "AM" = 113.6
    "AN" = -46.74117647058824
    "AO" = -1.001760563380282
    "AP" = -42.51176470588236
    "AQ" = -1.001760563380282
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range("${col}2").Value = $val
    $ws.Range("${col}3").Value = $val
}
